$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1076.6897
$ws.Range("J17").Value = 1076.6897
$ws.Range("L17").Value = 3230.0691
$ws.Range("N17").Value = -3566.0691
$ws.Range("H103").Value = 896.5
$ws.Range("I103").Value = 780
$ws.Range("J103").Value = 919.8
$ws.Range("K103").Value = 2340
$ws.Range("L103").Value = 2759.4
$ws.Range("M103").Value = -1754
$ws.Range("N103").Value = -3931.4
$ws.Range("H116").Value = 2808.4614
$ws.Range("I116").Value = 2463.9473
$ws.Range("J116").Value = 3743.5715
$ws.Range("K116").Value = 2463.9473
$ws.Range("L116").Value = 3743.5715
$ws.Range("M116").Value = 978.0527000000002
$ws.Range("N116").Value = -10627.5715
$ws.Range("H129").Value = 848.3108
$ws.Range("J129").Value = 889.98505
$ws.Range("L129").Value = 2669.95515
$ws.Range("N129").Value = -12669.95515
$ws.Range("H138").Value = 1055.0781
$ws.Range("I138").Value = 790.1064
$ws.Range("J138").Value = 1787.6471
$ws.Range("K138").Value = 2370.3192
$ws.Range("L138").Value = 5362.9413
$ws.Range("M138").Value = 2769.6808
$ws.Range("N138").Value = -15642.9413
$ws.Range("H141").Value = 862.5
$ws.Range("I141").Value = 850
$ws.Range("J141").Value = 1000
$ws.Range("K141").Value = 2550
$ws.Range("L141").Value = 3000
$ws.Range("M141").Value = 2630
$ws.Range("N141").Value = -13360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 100013
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("H4").Value = 500
$ws.Range("J4").Value = 500
$ws.Range("L4").Value = 500
$ws.Range("N4").Value = -732
$ws.Range("H5").Value = 200
$ws.Range("I5").Value = 200
$ws.Range("K5").Value = 200
$ws.Range("M5").Value = -88
$ws.Range("H15").Value = 13000
$ws.Range("J15").Value = 13000
$ws.Range("L15").Value = 13000
$ws.Range("N15").Value = -13700
$ws.Range("H74").Value = 995.6799999999999
$ws.Range("I74").Value = 813.9
$ws.Range("K74").Value = 813.9
$ws.Range("M74").Value = 60.10000000000002
$ws.Range("H75").Value = 29173
$ws.Range("J75").Value = 29173
$ws.Range("L75").Value = 29173
$ws.Range("N75").Value = -30921
$ws.Range("H77").Value = 995.6799999999999
$ws.Range("I77").Value = 813.9
$ws.Range("K77").Value = 4069.5
$ws.Range("M77").Value = 298.5
$ws.Range("H78").Value = 29173
$ws.Range("J78").Value = 29173
$ws.Range("L78").Value = 87519
$ws.Range("N78").Value = -96255
$ws.Range("H98").Value = 29117.334
$ws.Range("J98").Value = 29117.334
$ws.Range("L98").Value = 29117.334
$ws.Range("N98").Value = -35107.334
$ws.Range("H102").Value = 23811496
$ws.Range("I102").Value = 27779246
$ws.Range("K102").Value = 27779246
$ws.Range("M102").Value = -27777624
$ws.Range("H112").Value = 7268.5
$ws.Range("J112").Value = 7268.5
$ws.Range("L112").Value = 7268.5
$ws.Range("N112").Value = -10222.5
$ws.Range("H114").Value = 19405.092
$ws.Range("J114").Value = 19405.092
$ws.Range("L114").Value = 19405.092
$ws.Range("N114").Value = -28083.092
$ws.Range("H116").Value = 100013
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("H132").Value = 2851.875
$ws.Range("I132").Value = 2476.5
$ws.Range("K132").Value = 7429.5
$ws.Range("M132").Value = -4899.5
$ws.Range("M2").ClearContents()
$ws.Range("M116").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 100013
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("H4").Value = 200
$ws.Range("I4").Value = 200
$ws.Range("K4").Value = 200
$ws.Range("M4").Value = -85
$ws.Range("H86").Value = 2756
$ws.Range("I86").Value = 2862.4194
$ws.Range("J86").Value = 2389.4443
$ws.Range("K86").Value = 2862.4194
$ws.Range("L86").Value = 2389.4443
$ws.Range("M86").Value = -1739.4194
$ws.Range("N86").Value = -4635.4443
$ws.Range("H89").Value = 2756
$ws.Range("I89").Value = 2862.4194
$ws.Range("J89").Value = 2389.4443
$ws.Range("K89").Value = 14312.097
$ws.Range("L89").Value = 11947.2215
$ws.Range("M89").Value = -8696.097000000002
$ws.Range("N89").Value = -23179.2215
$ws.Range("H99").Value = 22728492
$ws.Range("I99").Value = 29412944
$ws.Range("J99").Value = 1362.4
$ws.Range("K99").Value = 29412944
$ws.Range("L99").Value = 1362.4
$ws.Range("M99").Value = -29411446
$ws.Range("N99").Value = -4358.4
$ws.Range("H105").Value = 83334740
$ws.Range("I105").Value = 100001150
$ws.Range("J105").Value = 2650
$ws.Range("K105").Value = 100001150
$ws.Range("L105").Value = 2650
$ws.Range("M105").Value = -99999403
$ws.Range("N105").Value = -6144
$ws.Range("H110").Value = 33379.8
$ws.Range("J110").Value = 33379.8
$ws.Range("L110").Value = 33379.8
$ws.Range("N110").Value = -41559.8
$ws.Range("M3").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2200.45
$ws.Range("I31").Value = 1223.3334
$ws.Range("J31").Value = 2999.9092
$ws.Range("K31").Value = 1223.3334
$ws.Range("L31").Value = 2999.9092
$ws.Range("M31").Value = -928.3334
$ws.Range("N31").Value = -3589.9092
$ws.Range("H34").Value = 2200.45
$ws.Range("I34").Value = 1223.3334
$ws.Range("J34").Value = 2999.9092
$ws.Range("K34").Value = 1223.3334
$ws.Range("L34").Value = 2999.9092
$ws.Range("M34").Value = -1021.3334
$ws.Range("N34").Value = -3403.9092
$ws.Range("H60").Value = 10080.5
$ws.Range("I60").Value = 2606.6
$ws.Range("J60").Value = 11860
$ws.Range("K60").Value = 2606.6
$ws.Range("L60").Value = 11860
$ws.Range("M60").Value = -2095.6
$ws.Range("N60").Value = -12882

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 2933.6667
$ws.Range("I58").Value = 905
$ws.Range("J58").Value = 3078.5715
$ws.Range("K58").Value = 2715
$ws.Range("L58").Value = 9235.7145
$ws.Range("M58").Value = -2587
$ws.Range("N58").Value = -9491.7145
$ws.Range("H131").Value = 17242696
$ws.Range("J131").Value = 1455.3137
$ws.Range("L131").Value = 4365.9411
$ws.Range("N131").Value = -14445.9411

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 6553.6665
$ws.Range("I43").Value = 3000
$ws.Range("J43").Value = 6997.875
$ws.Range("K43").Value = 3000
$ws.Range("L43").Value = 6997.875
$ws.Range("M43").Value = -2849
$ws.Range("N43").Value = -7299.875
$ws.Range("H80").Value = 4944.4287
$ws.Range("I80").Value = 2752.5
$ws.Range("J80").Value = 5821.2
$ws.Range("K80").Value = 2752.5
$ws.Range("L80").Value = 5821.2
$ws.Range("M80").Value = -1754.5
$ws.Range("N80").Value = -7817.2
$ws.Range("H83").Value = 4944.4287
$ws.Range("I83").Value = 2752.5
$ws.Range("J83").Value = 5821.2
$ws.Range("K83").Value = 13762.5
$ws.Range("L83").Value = 29106
$ws.Range("M83").Value = -8770.5
$ws.Range("N83").Value = -39090
$ws.Range("H103").Value = 10000
$ws.Range("J103").Value = 10000
$ws.Range("L103").Value = 10000
$ws.Range("N103").Value = -12344
$ws.Range("H111").Value = 20000
$ws.Range("J111").Value = 20000
$ws.Range("L111").Value = 20000
$ws.Range("N111").Value = -26134

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 27839.8
$ws.Range("J110").Value = 27299.75
$ws.Range("L110").Value = 27299.75
$ws.Range("N110").Value = -35479.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 7999.25
$ws.Range("J93").Value = 7999.25
$ws.Range("L93").Value = 7999.25
$ws.Range("N93").Value = -12991.25
$ws.Range("H97").Value = 20572
$ws.Range("J97").Value = 20572
$ws.Range("L97").Value = 20572
$ws.Range("N97").Value = -22554
